# PerformerConfig.xlsx update
# - Settings sheet: rename the QueueFolder value from "LazyFramework.Tests" to
#   "LazyFramework", and add two new rows describing a maintenance window.
# - Assets sheet: add two new LazyFramework assets for the IMAP server.
# - Minor page setup / selection tidy-up to mirror what Excel records when a
#   user finishes editing a sheet.

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("Settings")
$assets   = $wb.Worksheets.Item("Assets")

# 1) Settings!B4 used to read "LazyFramework.Tests" - fix it to "LazyFramework".
$settings.Range("B4").Value = "LazyFramework"

# 2) Assets: two new rows for the IMAP mailbox used by the framework.
#    (Name | Value | Folder | Description)
$assets.Cells.Item(10, 1).Value = "IMAP_Port"
$assets.Cells.Item(10, 2).Value = "IMAP_Port"
$assets.Cells.Item(10, 3).Value = "LazyFramework"
$assets.Cells.Item(10, 4).Value = "The port of the IMAP server."

$assets.Cells.Item(11, 1).Value = "IMAP_Server"
$assets.Cells.Item(11, 2).Value = "IMAP_Server"
$assets.Cells.Item(11, 3).Value = "LazyFramework"
$assets.Cells.Item(11, 4).Value = "The URL of the IMAP server."

# 3) Settings: two new rows describing a maintenance window.
#    Fill the Name column first, then the start/end descriptions, then the
#    (text-formatted) time values - this matches the order the strings were
#    originally typed in.
$settings.Cells.Item(13, 1).Value = "Maintenance_Start"
$settings.Cells.Item(14, 1).Value = "Maintenance_End"

$settings.Cells.Item(14, 3).Value = "The end of a maintenance window."
$settings.Cells.Item(13, 3).Value = "The start of a maintenance window."

$settings.Cells.Item(13, 2).NumberFormat = "@"
$settings.Cells.Item(13, 2).Value = "00:00:00"
$settings.Cells.Item(14, 2).NumberFormat = "@"
$settings.Cells.Item(14, 2).Value = "00:00:00"

# 4) Page setup on the Settings sheet (explicit portrait orientation).
$settings.PageSetup.Orientation = 1

# 5) Leave the selection where the author last clicked on each sheet.
$settings.Range("B14").Select()
$assets.Range("G16").Select()

Write-Host "PerformerConfig.xlsx updated"
